# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": bump Version/Date, insert a new "Jurisdiction" row
# right after "Contact" (pushes every row below it down by one).
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.1.1"
$meta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Shift rows 11..21 down to 12..22 (bottom-up so nothing is clobbered).
# Each shift copies both the values AND formatting of the source row onto
# the destination row, so the destination keeps reusing the same style
# that's already in the sheet instead of minting a brand-new cellXf.
for ($r = 21; $r -ge 11; $r--) {
    $nr = $r + 1
    $src = $meta.Range("A" + $r + ":B" + $r)
    $dst = $meta.Range("A" + $nr + ":B" + $nr)
    $dst.ClearContents()
    $src.Copy()
    $dst.PasteSpecial(-4163)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# Row 11 becomes the new "Jurisdiction" property (value left blank).
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").ClearContents()

# ---------------------------------------------------------------------
# Sheet "Concepts": lower-case / re-symbol the Code column.
# ---------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Range("B2").Value = "<4w"
$concepts.Range("B3").Value = "4-6w"
$concepts.Range("B4").Value = "7-12w"
$concepts.Range("B5").Value = "12w-6m"
$concepts.Range("B6").Value = ">6m"
$concepts.Range("B8").Value = "P"
